# Update Wnt5a-Lrp5 LR-pairs worksheet with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.977669497583861
$ws.Range("J2").Value = 0.977669497583861
$ws.Range("M2").Value = 10.25883033333333
$ws.Range("N2").Value = 30.776491
$ws.Range("O2").Value = 0.34684992242997
$ws.Range("P2").Value = 0.34684992242997
$ws.Range("Q2").Value = 77.09300689478167
$ws.Range("R2").Value = 693.837062053035
$ws.Range("S2").Value = 0.3391045893991099
$ws.Range("T2").Value = 0.3391045893991099
$ws.Range("I3").Value = 0.977669497583861
$ws.Range("J3").Value = 0.977669497583861
$ws.Range("O3").Value = 0.1682819529322607
$ws.Range("P3").Value = 0.1682819529322608
$ws.Range("S3").Value = 0.1645241323757143
$ws.Range("T3").Value = 0.1645241323757143
$ws.Range("I4").Value = 0.977669497583861
$ws.Range("J4").Value = 0.977669497583861
$ws.Range("M4").Value = 4.130648333333333
$ws.Range("N4").Value = 12.391945
$ws.Range("O4").Value = 0.1396567647041521
$ws.Range("P4").Value = 0.1396567647041521
$ws.Range("Q4").Value = 31.04097544209166
$ws.Range("R4").Value = 279.368778978825
$ws.Range("S4").Value = 0.1365381589824958
$ws.Range("T4").Value = 0.1365381589824959
$ws.Range("I5").Value = 0.977669497583861
$ws.Range("J5").Value = 0.977669497583861
$ws.Range("M5").Value = 2.760918333333333
$ws.Range("N5").Value = 8.282755
$ws.Range("O5").Value = 0.09334634443076846
$ws.Range("P5").Value = 0.09334634443076847
$ws.Range("Q5").Value = 20.74773528674166
$ws.Range("R5").Value = 186.729617580675
$ws.Range("S5").Value = 0.09126187366091944
$ws.Range("T5").Value = 0.09126187366091945
$ws.Range("I6").Value = 0.977669497583861
$ws.Range("J6").Value = 0.977669497583861
$ws.Range("M6").Value = 2.602884
$ws.Range("N6").Value = 7.808651999999999
$ws.Range("O6").Value = 0.08800322104565558
$ws.Range("P6").Value = 0.0880032210456556
$ws.Range("Q6").Value = 19.56013966878
$ws.Range("R6").Value = 176.04125701902
$ws.Range("S6").Value = 0.08603806490546756
$ws.Range("T6").Value = 0.08603806490546757
$ws.Range("I7").Value = 0.977669497583861
$ws.Range("J7").Value = 0.977669497583861
$ws.Range("M7").Value = 4.846564
$ws.Range("N7").Value = 14.539692
$ws.Range("O7").Value = 0.1638617944571931
$ws.Range("P7").Value = 0.1638617944571932
$ws.Range("Q7").Value = 36.42093491438
$ws.Range("R7").Value = 327.78841422942
$ws.Range("S7").Value = 0.1602026782601539
$ws.Range("T7").Value = 0.1602026782601539
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.171642
$ws.Range("H8").Value = 0.514926
$ws.Range("I8").Value = 0.02233050241613897
$ws.Range("J8").Value = 0.02233050241613898
$ws.Range("M8").Value = 10.25883033333333
$ws.Range("N8").Value = 30.776491
$ws.Range("O8").Value = 0.34684992242997
$ws.Range("P8").Value = 0.34684992242997
$ws.Range("Q8").Value = 1.760846156074
$ws.Range("R8").Value = 15.847615404666
$ws.Range("S8").Value = 0.00774533303086006
$ws.Range("T8").Value = 0.007745333030860061
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.171642
$ws.Range("H9").Value = 0.514926
$ws.Range("I9").Value = 0.02233050241613897
$ws.Range("J9").Value = 0.02233050241613898
$ws.Range("O9").Value = 0.1682819529322607
$ws.Range("P9").Value = 0.1682819529322608
$ws.Range("Q9").Value = 0.8543136693860001
$ws.Range("R9").Value = 7.688823024474001
$ws.Range("S9").Value = 0.003757820556546433
$ws.Range("T9").Value = 0.003757820556546434
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.171642
$ws.Range("H10").Value = 0.514926
$ws.Range("I10").Value = 0.02233050241613897
$ws.Range("J10").Value = 0.02233050241613898
$ws.Range("M10").Value = 4.130648333333333
$ws.Range("N10").Value = 12.391945
$ws.Range("O10").Value = 0.1396567647041521
$ws.Range("P10").Value = 0.1396567647041521
$ws.Range("Q10").Value = 0.7089927412299999
$ws.Range("R10").Value = 6.380934671069999
$ws.Range("S10").Value = 0.00311860572165622
$ws.Range("T10").Value = 0.00311860572165622
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.171642
$ws.Range("H11").Value = 0.514926
$ws.Range("I11").Value = 0.02233050241613897
$ws.Range("J11").Value = 0.02233050241613898
$ws.Range("M11").Value = 2.760918333333333
$ws.Range("N11").Value = 8.282755
$ws.Range("O11").Value = 0.09334634443076846
$ws.Range("P11").Value = 0.09334634443076847
$ws.Range("Q11").Value = 0.4738895445699999
$ws.Range("R11").Value = 4.26500590113
$ws.Range("S11").Value = 0.002084470769849016
$ws.Range("T11").Value = 0.002084470769849016
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.171642
$ws.Range("H12").Value = 0.514926
$ws.Range("I12").Value = 0.02233050241613897
$ws.Range("J12").Value = 0.02233050241613898
$ws.Range("M12").Value = 2.602884
$ws.Range("N12").Value = 7.808651999999999
$ws.Range("O12").Value = 0.08800322104565558
$ws.Range("P12").Value = 0.0880032210456556
$ws.Range("Q12").Value = 0.4467642155279999
$ws.Range("R12").Value = 4.020877939751999
$ws.Range("S12").Value = 0.001965156140188024
$ws.Range("T12").Value = 0.001965156140188025
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.171642
$ws.Range("H13").Value = 0.514926
$ws.Range("I13").Value = 0.02233050241613897
$ws.Range("J13").Value = 0.02233050241613898
$ws.Range("M13").Value = 4.846564
$ws.Range("N13").Value = 14.539692
$ws.Range("O13").Value = 0.1638617944571931
$ws.Range("P13").Value = 0.1638617944571932
$ws.Range("Q13").Value = 0.8318739380879999
$ws.Range("R13").Value = 7.486865442792
$ws.Range("S13").Value = 0.003659116197039219
$ws.Range("T13").Value = 0.00365911619703922

Write-Output "Updated cells with new TPM values."
